# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Mapping of row -> new value (same change applied to both sheets).

$wb = $excel.ActiveWorkbook

$changes = @{
    8  = 78
    10 = 45
    11 = 20
    12 = 560
    13 = 26
    14 = 293
    15 = 25
    16 = 354
    21 = 89
    22 = 884
    23 = 1387
    24 = 294
    25 = 319
    27 = 68
    29 = 36
    31 = 207
    32 = 239
    33 = 264
    34 = 1603
    35 = 49
    40 = 3534
    41 = 415
    42 = 185
    43 = 897
    45 = 62
    46 = 41
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Cells.Item($row, 6).Value = $changes[$row]
    }
}
